$wb = $excel.ActiveWorkbook

# The two sheets "展览" and "全部类型" contain identical data tables and
# both need the same "想去人数" (want-to-go count) updates applied.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 2960
    $ws.Range("F5").Value = 6721
    $ws.Range("F6").Value = 1697
    $ws.Range("F7").Value = 21
}
